$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("B5").Value = -0.1818718991761584
$ws.Range("C5").Value = 2.0051076125427789
$ws.Range("D5").Value = 614895744
$ws.Range("E5").Value = 2464063
$ws.Range("G5").Value = -0.19156646969580171
$ws.Range("H5").Value = 1.862128866185955
$ws.Range("I5").Value = 413674742
$ws.Range("J5").Value = 1652722

# Row 13 updates
$ws.Range("B13").Value = -0.29549899240109923
$ws.Range("C13").Value = 1.728488423210542
$ws.Range("D13").Value = 867179030
$ws.Range("E13").Value = 2854490
$ws.Range("G13").Value = -0.27674912109777727
$ws.Range("H13").Value = 1.843608639394998
$ws.Range("I13").Value = 1060203761
$ws.Range("J13").Value = 3638818
